$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vancouver CTs Matched to CSDs -")

# Rows where CSDSimple (column C) holds the District Municipality flavour of
# "North Vancouver" (CSDTYPE = DM) -> disambiguate to "North Vancovuer District"
$dmRows = @(7, 25, 45, 59, 75, 76, 86, 154, 162, 167, 221, 224, 226, 246, 287, 310, 341, 370)

# Rows where CSDSimple (column C) holds the City flavour of
# "North Vancouver" (CSDTYPE = CY) -> disambiguate to "North Vancouver City"
$cyRows = @(13, 27, 141, 155, 178, 207, 216, 401, 420)

foreach ($r in $dmRows) {
    $ws.Cells.Item($r, 3).Value = "North Vancovuer District"
}

foreach ($r in $cyRows) {
    $ws.Cells.Item($r, 3).Value = "North Vancouver City"
}

# Reproduce the selection change recorded in the saved file
$ws.Range("D7").Select()

# Turn on the AutoFilter over the data range, which is what produces the
# hidden _xlnm._FilterDatabase defined name seen in the diff
$ws.Range("A1:F458").AutoFilter() | Out-Null

# AutoFilter registers a sheet-scoped, hidden defined name
# (_xlnm._FilterDatabase) pointing at the filtered range
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "='Vancouver CTs Matched to CSDs -'!`$A`$1:`$F`$458")
$filterName.Visible = $false

$wb.Save()
